$wb = $excel.ActiveWorkbook

# --- Sheet "Vanilla" (sheet1): update B4 value and active selection ---
$ws1 = $wb.Worksheets.Item("Vanilla")
$ws1.Range("B4").Value = 100
$ws1.Activate()
$ws1.Range("B5").Select()

# --- Sheet "P8_Split_P6.2_only" (sheet2): insert new row and update values ---
$ws2 = $wb.Worksheets.Item("P8_Split_P6.2_only")

# Insert a new row before row 6, shifting rows 6-10 down to 7-11
$ws2.Rows.Item(6).Insert()

# Set the new row 6 contents: "p6_infants" in A6, leave B6 empty
$ws2.Range("A6").Value = "p6_infants"

# Update B4 value (same as sheet1 change)
$ws2.Range("B4").Value = 100

$ws2.Activate()
$ws2.Range("A6").Select()
